$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.56"
$ws.Range("E2").Value = "'4.85%"
$ws.Range("D3").Value = "'27.26"
$ws.Range("E3").Value = "'-4.01%"
$ws.Range("D4").Value = "'5.215"
$ws.Range("E4").Value = "'-1.40%"
$ws.Range("D5").Value = "'0.05920"
$ws.Range("E5").Value = "'3.66%"
$ws.Range("E6").Value = "'0.63%"
$ws.Range("D7").Value = "'0.8641"
$ws.Range("E7").Value = "'1.12%"
$ws.Range("E8").Value = "'13.37%"
$ws.Range("D9").Value = "'0.1414"
$ws.Range("E9").Value = "'1.56%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07181"
$ws.Range("E10").Value = "'1.24%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03113"
$ws.Range("E11").Value = "'-1.54%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09234"
$ws.Range("E12").Value = "'0.05%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001539"
$ws.Range("E13").Value = "'0.79%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006077"
$ws.Range("E14").Value = "'1.29%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005667"
$ws.Range("E15").Value = "'-3.64%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.495"
$ws.Range("E16").Value = "'-0.06%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.268"
$ws.Range("E17").Value = "'1.78%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.192"
$ws.Range("E18").Value = "'0.93%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3146"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03566"
$ws.Range("E20").Value = "'6.81%"
$ws.Range("E21").Value = "'-0.05%"
$ws.Range("D22").Value = "'3.517"
$ws.Range("E22").Value = "'1.06%"
$ws.Range("D23").Value = "'0.04182"
$ws.Range("E23").Value = "'2.92%"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'0.06%"
$ws.Range("E26").Value = "'8.63%"
$ws.Range("E27").Value = "'-0.07%"
$ws.Range("E28").Value = "'2.68%"
$ws.Range("D40").Value = "'0.03812"
$ws.Range("E40").Value = "'0.46%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1104"
$ws.Range("E41").Value = "'3.46%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003995"
$ws.Range("E42").Value = "'-22.49%"
$ws.Range("D43").Value = "'0.002358"
$ws.Range("E43").Value = "'-2.95%"
$ws.Range("D44").Value = "'0.01082"
$ws.Range("E44").Value = "'14.70%"
$ws.Range("D45").Value = "'0.00005425"
$ws.Range("E45").Value = "'2.71%"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("E47").Value = "'22.39%"
$ws.Range("D48").Value = "'0.002234"
$ws.Range("E48").Value = "'-1.05%"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E50").Value = "'-0.05%"
